$wb = $excel.ActiveWorkbook

# --- Sheet "2EXT04_DNA": add content examples + extend the annotation table ---
$ws1 = $wb.Worksheets.Item("2EXT04_DNA")

# Extend the Swate annotation table from A1:S2 to A1:S5 so the new example rows
# become part of the table (and the autofilter range grows with it).
$lo = $ws1.ListObjects.Item(1)
$lo.Resize($ws1.Range("A1:S5"))

# Row 2 already carries the formatted (unit) cells E2:R2 - copy that formatting
# down to the new rows 3-5 before filling in the example values.
$ws1.Range("E2:R2").Copy($ws1.Range("E3:R5"))

# Row 2: DNA/RNA extraction example using QIAGEN RNEasy
$ws1.Range("B2").Value = "RNA (Transcriptomics)"
$ws1.Range("I2").Value = "QIAGEN RNEasy"
$ws1.Range("L2").Value = "QIAGEN RNEasy Buffer 2"
$ws1.Range("O2").Value = 200

# Row 3: DNA extraction example using PCI method
$ws1.Range("B3").Value = "DNA (Genomics)"
$ws1.Range("I3").Value = "PCI method"
$ws1.Range("L3").Value = "phenol:chloroform:isopropanol"

# Row 4 & 5: additional bio-entity examples
$ws1.Range("B4").Value = "Metabolites"
$ws1.Range("B5").Value = "Protein"

# --- Sheet "SwateTemplateMetadata": bump template version ---
$ws2 = $wb.Worksheets.Item("SwateTemplateMetadata")
$ws2.Range("B3").Value = "1.1.4"
